$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

$ws.Range("E2").Value = '94.20-1-00'
$ws.Range("F2").Value = 'S'
$ws.Range("G2").Value = 'OUTRAS ATIVIDADES DE SERVIÇOS'
$ws.Range("H2").Value = 94
$ws.Range("I2").Value = 'ATIVIDADES DE ORGANIZAÇÕES ASSOCIATIVAS'
$ws.Range("E3").Value = '47.11-3-02'
$ws.Range("F3").Value = 'G'
$ws.Range("G3").Value = 'COMÉRCIO; REPARAÇÃO DE VEÍCULOS AUTOMOTORES E MOTOCICLETAS'
$ws.Range("H3").Value = 47
$ws.Range("I3").Value = 'COMÉRCIO VAREJISTA'
$ws.Range("J3").Value = 'COMÉRCIO VAREJISTA DE MERCADORIAS EM GERAL, COM PREDOMINÂNCIA DE PRODUTOS ALIMENTÍCIOS - SUPERMERCADOS'
$ws.Range("B6").Value = 'DANIEL BORGONOVO THIVES'
$ws.Range("C6").Value = 'JOINVILLE '
$ws.Range("D6").Value = ' SC'
$ws.Range("E6").Value = '56.11-2-01'
$ws.Range("F6").Value = 'I'
$ws.Range("G6").Value = 'ALOJAMENTO E ALIMENTAÇÃO'
$ws.Range("H6").Value = 56
$ws.Range("I6").Value = 'ALIMENTAÇÃO'
$ws.Range("J6").Value = 'RESTAURANTES E SIMILARES'
$ws.Range("B7").Value = 'HP PNEUS E SERVICOS AUTOMOTIVOS LTDA'
$ws.Range("C7").Value = 'BRASILIA '
$ws.Range("D7").Value = ' DF'
$ws.Range("B8").Value = 'M.S.M.- MINAS SERVICOS E MONTAGENS LTDA'
$ws.Range("C8").Value = 'ITAJUBA '
$ws.Range("D8").Value = ' MG'
$ws.Range("E8").Value = '26.32-9-00'
$ws.Range("F8").Value = 'C'
$ws.Range("G8").Value = 'INDÚSTRIAS DE TRANSFORMAÇÃO'
$ws.Range("H8").Value = 26
$ws.Range("I8").Value = 'FABRICAÇÃO DE EQUIPAMENTOS DE INFORMÁTICA, PRODUTOS ELETRÔNICOS E ÓPTICOS'
$ws.Range("J8").Value = 'FABRICAÇÃO DE APARELHOS TELEFÔNICOS E DE OUTROS EQUIPAMENTOS DE COMUNICAÇÃO, PEÇAS E ACESSÓRIOS'
$ws.Range("B9").Value = 'IRMANDADE DA SANTA CASA DE CARIDADE DE MACHADO'
$ws.Range("C9").Value = 'MACHADO '
$ws.Range("D9").Value = ' MG'
$ws.Range("B10").Value = 'SINDICATO DOS TRABALHADORES NAS EMPRESAS DE ENERGIA DO RIO DE JANEIRO E REGIAO - SINTERGIA/RJ'
$ws.Range("C10").Value = 'RIO DE JANEIRO '
$ws.Range("D10").Value = ' RJ'
$ws.Range("E10").Value = '94.20-1-00'
$ws.Range("F10").Value = 'S'
$ws.Range("G10").Value = 'OUTRAS ATIVIDADES DE SERVIÇOS'
$ws.Range("H10").Value = 94
$ws.Range("I10").Value = 'ATIVIDADES DE ORGANIZAÇÕES ASSOCIATIVAS'
$ws.Range("J10").Value = 'ATIVIDADES DE ORGANIZAÇÕES SINDICAIS'
$ws.Range("B11").Value = 'NIDEC GLOBAL APPLIANCE BRASIL LTDA'
$ws.Range("C11").Value = 'ITAIOPOLIS '
$ws.Range("D11").Value = ' SC'
$ws.Range("B12").Value = 'NIDEC GLOBAL APPLIANCE BRASIL LTDA'
$ws.Range("C12").Value = 'JOINVILLE '
$ws.Range("D12").Value = ' SC'
$ws.Range("E12").Value = '24.51-2-00'
$ws.Range("F12").Value = 'C'
$ws.Range("G12").Value = 'INDÚSTRIAS DE TRANSFORMAÇÃO'
$ws.Range("H12").Value = 24
$ws.Range("I12").Value = 'METALURGIA'
$ws.Range("J12").Value = 'FUNDIÇÃO DE FERRO E AÇO'
$ws.Range("B13").Value = 'NIDEC GLOBAL APPLIANCE BRASIL LTDA'
$ws.Range("C13").Value = 'JOINVILLE '
$ws.Range("D13").Value = ' SC'
$ws.Range("B14").Value = 'TUBOCANO ARTEFATOS DE CIMENTO LTDA'
$ws.Range("C14").Value = 'GRAVATAI '
$ws.Range("D14").Value = ' RS'
$ws.Range("B15").Value = 'UNIVERSAL LEAF TABACOS LTDA'
$ws.Range("C15").Value = 'SANTA CRUZ DO SUL '
$ws.Range("D15").Value = ' RS'
$ws.Range("B16").Value = 'PELICAN TEXTIL LTDA'
$ws.Range("C16").Value = 'SANTA ISABEL '
$ws.Range("D16").Value = ' SP'
$ws.Range("B17").Value = 'HERCILIO RESTAURANTE E CAFE LTDA'
$ws.Range("C17").Value = 'FLORIANOPOLIS '
$ws.Range("D17").Value = ' SC'
$ws.Range("B18").Value = 'FINSOL SOCIEDADE DE CREDITO AO MICROEMPREENDEDOR E A EMPRESA DE PEQUENO PORTE S/A'
$ws.Range("C18").Value = 'RECIFE '
$ws.Range("D18").Value = ' PE'
$ws.Range("B19").Value = 'RJS CONSTRUCAO CIVIL LTDA'
$ws.Range("C19").Value = 'FRANCISCO MORATO '
$ws.Range("D19").Value = ' SP'
$ws.Range("E19").Value = '41.20-4-00'
$ws.Range("F19").Value = 'F'
$ws.Range("G19").Value = 'CONSTRUÇÃO'
$ws.Range("H19").Value = 41
$ws.Range("I19").Value = 'CONSTRUÇÃO DE EDIFÍCIOS'
$ws.Range("J19").Value = 'CONSTRUÇÃO DE EDIFÍCIOS'
$ws.Range("B20").Value = 'CESARI ENGENHARIA E CONSTRUCAO LTDA'
$ws.Range("C20").Value = 'BARRETOS '
$ws.Range("D20").Value = ' SP'
$ws.Range("B21").Value = 'PECUARIA FERNANDO LTDA'
$ws.Range("C21").Value = 'MEDIANEIRA '
$ws.Range("D21").Value = ' PR'
$ws.Range("B22").Value = 'B-GREEN GESTAO AMBIENTAL LTDA'
$ws.Range("C22").Value = 'RECIFE '
$ws.Range("D22").Value = ' PE'
$ws.Range("B23").Value = 'L & R RESTAURANTE EIRELI'
$ws.Range("C23").Value = 'CAMPO GRANDE '
$ws.Range("D23").Value = ' MS'
$ws.Range("B24").Value = 'TIARAJU ENGENHARIA LTDA'
$ws.Range("C24").Value = 'PASSO FUNDO '
$ws.Range("D24").Value = ' RS'
$ws.Range("B25").Value = 'SIND ESTAB ENSINO NO EST DO RIO DE JANEIRO'
$ws.Range("C25").Value = 'NITEROI '
$ws.Range("D25").Value = ' RJ'
$ws.Range("E25").Value = '94.20-1-00'
$ws.Range("F25").Value = 'S'
$ws.Range("G25").Value = 'OUTRAS ATIVIDADES DE SERVIÇOS'
$ws.Range("H25").Value = 94
$ws.Range("I25").Value = 'ATIVIDADES DE ORGANIZAÇÕES ASSOCIATIVAS'
$ws.Range("J25").Value = 'ATIVIDADES DE ORGANIZAÇÕES SINDICAIS'
$ws.Range("B26").Value = 'MATEUBRAS - COMERCIO DE MATERIAIS PARA CONSTRUCAO LTDA'
$ws.Range("C26").Value = 'SAO PAULO '
$ws.Range("D26").Value = ' SP'
$ws.Range("B27").Value = 'HTS ADMINISTRADORA DE HOTEIS LTDA.'
$ws.Range("C27").Value = 'MANAUS '
$ws.Range("D27").Value = ' AM'
$ws.Range("B28").Value = 'ASSOCIACAO SULINA DE CREDITO E ASSISTENCIA RURAL'
$ws.Range("C28").Value = 'PORTO ALEGRE '
$ws.Range("D28").Value = ' RS'
$ws.Range("B29").Value = 'NUTRIPURA NUTRICAO ANIMAL LTDA'
$ws.Range("C29").Value = 'RONDONOPOLIS '
$ws.Range("D29").Value = ' MT'
$ws.Range("B30").Value = 'SINDICATO DOS EMPRESARIOS LOTERICOS DO ESTADO DE GOIAS'
$ws.Range("C30").Value = 'GOIANIA '
$ws.Range("D30").Value = ' GO'
$ws.Range("E30").Value = '94.20-1-00'
$ws.Range("F30").Value = 'S'
$ws.Range("G30").Value = 'OUTRAS ATIVIDADES DE SERVIÇOS'
$ws.Range("H30").Value = 94
$ws.Range("I30").Value = 'ATIVIDADES DE ORGANIZAÇÕES ASSOCIATIVAS'
$ws.Range("J30").Value = 'ATIVIDADES DE ORGANIZAÇÕES SINDICAIS'
$ws.Range("B31").Value = 'CHT QUIMIPEL BRAZIL QUIMICA LTDA.'
$ws.Range("C31").Value = 'PIRACAIA '
$ws.Range("D31").Value = ' SP'
$ws.Range("B32").Value = 'CSVIVA SOLUCOES TECNICAS LTDA'
$ws.Range("C32").Value = 'SAO CARLOS '
$ws.Range("D32").Value = ' SP'
$ws.Range("B33").Value = 'CONSERV LOCADORA E TRANSPORTADORA LTDA'
$ws.Range("C33").Value = 'JACOBINA '
$ws.Range("D33").Value = ' BA'
$ws.Range("B34").Value = 'CONSERV LOCADORA E TRANSPORTADORA LTDA'
$ws.Range("C34").Value = 'MARACAS '
$ws.Range("D34").Value = ' BA'
